$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.768000000000001
$ws.Range("A9").Value = -21.658
$ws.Range("D12").Value = -7.456
$ws.Range("A13").Value = -22.219
$ws.Range("D14").Value = -7.782999999999999
$ws.Range("A16").Value = -22.027
$ws.Range("A18").Value = -22.086
$ws.Range("D19").Value = -7.851999999999999
$ws.Range("A20").Value = -20.1
$ws.Range("A26").Value = -21.235
$ws.Range("D26").Value = -8.122999999999999
$ws.Range("A27").Value = -21.705
$ws.Range("D27").Value = -8.599
$ws.Range("A29").Value = -21.347
$ws.Range("D29").Value = -7.329000000000001
$ws.Range("A35").Value = -19.823
$ws.Range("A36").Value = -20.652
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.722
$ws.Range("A45").Value = -21.573
$ws.Range("D47").Value = -7.502
$ws.Range("D51").Value = -8.638000000000002
$ws.Range("D52").Value = -8.1
$ws.Range("A55").Value = -22.182
$ws.Range("D55").Value = -8.141000000000002
$ws.Range("A57").Value = -22.258
$ws.Range("A69").Value = -21.579
$ws.Range("D69").Value = -7.528999999999999
$ws.Range("D70").Value = -7.456
$ws.Range("A76").Value = -20.047
$ws.Range("D76").Value = -7.742999999999999
$ws.Range("A78").Value = -20.182
$ws.Range("A82").Value = -22.018
$ws.Range("A83").Value = -21.853
$ws.Range("D83").Value = -8.462999999999999
$ws.Range("A93").Value = -21.559
$ws.Range("D94").Value = -7.456
$ws.Range("A97").Value = -22.044
$ws.Range("D100").Value = -8.294
$ws.Range("D102").Value = -7.906999999999999
